# shop 데이터 notice, intro NULL로 수정
# Fill the blank "notice" (P) / "intro" (Q) cells with the literal text "Null",
# and change the one existing "NULL" entry (P22) to "Null" as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose P (notice) column needs to be set to "Null" (was blank)
$pRows = @(3, 10, 11, 12, 14, 15, 16, 17, 19)
foreach ($r in $pRows) {
    $ws.Cells.Item($r, 16).Value = "Null"
}

# Rows whose Q (intro) column needs to be set to "Null" (was blank).
# These cells also pick up the s=19 border style (thin left edge) instead of
# the previous s=23 (medium left edge) - matches what column P already uses.
$qRows = @(4, 5, 6, 10, 11, 12, 13, 14, 15, 16, 17, 19)
foreach ($r in $qRows) {
    $cell = $ws.Cells.Item($r, 17)
    $cell.Value = "Null"
    $cell.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
}

# P22 previously held the shared string "NULL" - normalise it to "Null" too.
$ws.Cells.Item(22, 16).Value = "Null"

# Restore the view state: selection moves from Q13 to Q19 and the window is
# scrolled so L7 is the top-left visible cell.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("Q19").Select()
